# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Vega Monumental Concepción - Alcachofa"
# as row 103, pushing the existing rows 103-127 down to 104-128.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row above the current row 103 (shifts 103:127 -> 104:128)
$ws.Rows("103:103").Insert()

# Populate the newly inserted row with the new record's values
$ws.Range("A103").Value = 11
$ws.Range("B103").Value = "Vega Monumental Concepción"
$ws.Range("C103").Value = "Bíobío"
$ws.Range("D103").Value = 45211
$ws.Range("E103").Value = 8
$ws.Range("F103").Value = 100112013
$ws.Range("G103").Value = "Alcachofa"
$ws.Range("H103").Value = "Española"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 60
$ws.Range("K103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = 15000
$ws.Range("N103").Value = "$/caja 30 unidades"
$ws.Range("O103").Value = "Provincia de Limarí"
$ws.Range("P103").Value = 500
$ws.Range("Q103").Value = 30
$ws.Range("R103").Value = "Hortaliza"

# Make sure the date cell keeps the date/time number format used by the rest of column D
$ws.Range("D103").NumberFormat = $ws.Range("D104").NumberFormat()
